$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q3" sheet as the 2nd tab (after "总计"),
#    pushing the existing "2022-Q1" / "2021-Q4" sheets one slot to the
#    right.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item(2))
$newSheet.Name = "2022-Q3"

# Re-fetch the (now 3rd) "2022-Q1" sheet to borrow its header/row
# formatting (bold, centered, bordered header row + index-column style)
# for the freshly inserted, blank sheet.
$q1 = $wb.Worksheets.Item(3)
$q1.Range("A1:H2").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Fill in the "2022-Q3" sheet contents.
# ---------------------------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'004536"
$newSheet.Range("C2").Value = "嘉实中小企业量化活力灵活配置混合"
$newSheet.Range("D2").Value = "'0.21"
$newSheet.Range("E2").Value = "'94.75"
$newSheet.Range("F2").Value = "'3.94"
$newSheet.Range("G2").Value = "'0.0083"
$newSheet.Range("H2").Value = 9

# ---------------------------------------------------------------------
# 3. Update the "总计" (summary) sheet: insert a new first data row for
#    2022-Q3 and shift the existing 2022-Q1 / 2021-Q4 rows down by one.
# ---------------------------------------------------------------------
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.01

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q1"
$total.Range("C3").Value = 9
$total.Range("D3").Value = 1.09

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q4"
$total.Range("C4").Value = 5
$total.Range("D4").Value = 0.3

# A2/A3 already carry the index-column style from the original sheet;
# copy that same formatting onto the newly-created A4 cell.
$total.Range("A3").Copy()
$total.Range("A4").PasteSpecial(-4122)
